# Re-analyzed carbohydrate esters & amides: split each combined category
# ("carboEster", "amide") back into its two original sub-categories and
# re-run the factorial ANOVA / Tukey results for them on the
# "litterChemistry" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("litterChemistry")

# --- Insert two new rows after the existing "carboEster" row (row 4) --------
# This pushes lipid/alkane/amide down from rows 5-7 to rows 7-9.
$ws.Range("A5:H6").EntireRow.Insert()

$ws.Range("A5").Value = "carboEster1"
$ws.Range("B5").Value = "*"
$ws.Range("C5").Value = "***"
$ws.Range("D5").Value = "o"
$ws.Range("E5").Value = "o"
$ws.Range("F5").Value = "o"
$ws.Range("G5").Value = "o"
$ws.Range("H5").Value = "o"

$ws.Range("A6").Value = "carboEster2"
$ws.Range("B6").Value = "*"
$ws.Range("C6").Value = "***"
$ws.Range("D6").Value = "***"
$ws.Range("E6").Value = "o"
$ws.Range("F6").Value = "o"
$ws.Range("G6").Value = "**"
$ws.Range("H6").Value = "o"

# --- Insert two new rows after the (now shifted) "amide" row (row 9) -------
# Before this insert: row7=lipid, row8=alkane, row9=amide.
$ws.Range("A10:H11").EntireRow.Insert()

$ws.Range("A10").Value = "amide1"
$ws.Range("B10").Value = "*"
$ws.Range("C10").Value = "*"
$ws.Range("D10").Value = "***"
$ws.Range("E10").Value = "o"
$ws.Range("F10").Value = "o"
$ws.Range("G10").Value = "o"
$ws.Range("H10").Value = "o"

$ws.Range("A11").Value = "amide2"
$ws.Range("B11").Value = "o"
$ws.Range("C11").Value = "o"
$ws.Range("D11").Value = "***"
$ws.Range("E11").Value = "o"
$ws.Range("F11").Value = "o"
$ws.Range("G11").Value = "o"
$ws.Range("H11").Value = "o"

# --- Apply a thin box border around every cell, bold on the header row -----
$header = $ws.Range("A1:H1")
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

$body = $ws.Range("A2:H11")
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2

# --- Selection / active cell, matching the saved view ----------------------
$ws.Range("D5").Select()
